$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = [double]"0.2466370749492557"
$ws.Range("C3").Value = [double]"0.009846537327280743"
$ws.Range("C4").Value = [double]"0.00072850246583498"
$ws.Range("C5").Value = [double]"0.02816842647723668"
$ws.Range("C6").Value = [double]"1.708316993799091e-18"
$ws.Range("C7").Value = [double]"0.000774521279597736"
$ws.Range("C8").Value = [double]"0.02311043898291583"
$ws.Range("C9").Value = [double]"0.1101516327223209"
$ws.Range("C10").Value = [double]"0.0002008355530203714"
$ws.Range("C11").Value = [double]"0.02204372345666923"
$ws.Range("C12").Value = [double]"3.01153745330273e-17"
$ws.Range("C13").Value = [double]"7.692741145799391e-05"
$ws.Range("C14").Value = [double]"0.326582795141821"
$ws.Range("C15").Value = [double]"0.0002056051372634622"
$ws.Range("C16").Value = [double]"0.01974258369999374"
$ws.Range("C17").Value = [double]"0.2117303953953317"
